# Site update for FAQ
# Update the recorded run Start time / End time / Duration values.

$d = $word.ActiveDocument

$d.Content.Find.Execute("Start time: 2017-12-27 19:15:05", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start time: 2018-01-31 13:04:14", 2)

$d.Content.Find.Execute("End time: 2017-12-27 19:41:36", $true, $false, $false, $false, $false,
                         $true, 1, $false, "End time: 2018-01-31 13:29:25", 2)

$d.Content.Find.Execute("Duration: 26.52 mins", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Duration: 25.19 mins", 2)
